$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data range before writing the refreshed programme (rows 2-27)
$ws.Range("A2:J27").ClearContents()

# Row 2
$ws.Cells.Item(2,1).Value = "'2025-02-04"
$ws.Cells.Item(2,2).Value = "'21:00"
$ws.Cells.Item(2,3).Value = "La vie après Siham"
$ws.Cells.Item(2,4).Value = "VO"
$ws.Cells.Item(2,5).Value = "'"
$ws.Cells.Item(2,6).Value = "Namir Abdel Messeeh"
$ws.Cells.Item(2,7).Value = "Invité ADRC Cédric Lépine - Partenariat ADRC"
$ws.Cells.Item(2,8).Value = "'"
$ws.Cells.Item(2,9).Value = "'"
$ws.Cells.Item(2,10).Value = "'"

# Row 3
$ws.Cells.Item(3,1).Value = "'2026-02-05"
$ws.Cells.Item(3,2).Value = "'21:00"
$ws.Cells.Item(3,3).Value = "Baise-en-ville"
$ws.Cells.Item(3,4).Value = "VF"
$ws.Cells.Item(3,5).Value = "'"
$ws.Cells.Item(3,6).Value = "Martin Jauvat"
$ws.Cells.Item(3,7).Value = "Invité ADRC Chloé Caye"
$ws.Cells.Item(3,8).Value = "'"
$ws.Cells.Item(3,9).Value = "'"
$ws.Cells.Item(3,10).Value = "'"

# Row 4
$ws.Cells.Item(4,1).Value = "'2026-02-06"
$ws.Cells.Item(4,2).Value = "'15:00"
$ws.Cells.Item(4,3).Value = "La pire mère au monde"
$ws.Cells.Item(4,4).Value = "VF"
$ws.Cells.Item(4,5).Value = "'"
$ws.Cells.Item(4,6).Value = "Pierre Mazingarbe"
$ws.Cells.Item(4,7).Value = "'"
$ws.Cells.Item(4,8).Value = "Séance EPHAD HANDI"
$ws.Cells.Item(4,9).Value = "'"
$ws.Cells.Item(4,10).Value = "'"

# Row 5
$ws.Cells.Item(5,1).Value = "'2026-02-06"
$ws.Cells.Item(5,2).Value = "'21:00"
$ws.Cells.Item(5,3).Value = "Victor comme tout le monde"
$ws.Cells.Item(5,4).Value = "VF"
$ws.Cells.Item(5,5).Value = "'"
$ws.Cells.Item(5,6).Value = "'"
$ws.Cells.Item(5,7).Value = "'"
$ws.Cells.Item(5,8).Value = "'"
$ws.Cells.Item(5,9).Value = "'"
$ws.Cells.Item(5,10).Value = "Avant-Première en Compétition"

# Row 6
$ws.Cells.Item(6,1).Value = "'2026-02-07"
$ws.Cells.Item(6,2).Value = "'15:00"
$ws.Cells.Item(6,3).Value = "Le grand Phuket"
$ws.Cells.Item(6,4).Value = "VO"
$ws.Cells.Item(6,5).Value = "'"
$ws.Cells.Item(6,6).Value = "Liu Yaonan"
$ws.Cells.Item(6,7).Value = "'"
$ws.Cells.Item(6,8).Value = "Ciné Jeunes du samedi 18h"
$ws.Cells.Item(6,9).Value = "'"
$ws.Cells.Item(6,10).Value = "Avant-Première en Compétition"

# Row 7
$ws.Cells.Item(7,1).Value = "'2026-02-07"
$ws.Cells.Item(7,2).Value = "'17:00"
$ws.Cells.Item(7,3).Value = "En route"
$ws.Cells.Item(7,4).Value = "VF"
$ws.Cells.Item(7,5).Value = "'"
$ws.Cells.Item(7,6).Value = "Alexei Mironov"
$ws.Cells.Item(7,7).Value = "'"
$ws.Cells.Item(7,8).Value = "Ciné goûter JP"
$ws.Cells.Item(7,9).Value = "'"
$ws.Cells.Item(7,10).Value = "Avant-Première"

# Row 8
$ws.Cells.Item(8,1).Value = "'2026-02-07"
$ws.Cells.Item(8,2).Value = "'18:00"
$ws.Cells.Item(8,3).Value = "Urchin"
$ws.Cells.Item(8,4).Value = "VO"
$ws.Cells.Item(8,5).Value = "'"
$ws.Cells.Item(8,6).Value = "Harris Dickinson"
$ws.Cells.Item(8,7).Value = "'"
$ws.Cells.Item(8,8).Value = "'"
$ws.Cells.Item(8,9).Value = "'"
$ws.Cells.Item(8,10).Value = "Avant-Première en Compétition"

# Row 9
$ws.Cells.Item(9,1).Value = "'2026-02-07"
$ws.Cells.Item(9,2).Value = "'21:15"
$ws.Cells.Item(9,3).Value = "Noise"
$ws.Cells.Item(9,4).Value = "VO"
$ws.Cells.Item(9,5).Value = "'"
$ws.Cells.Item(9,6).Value = "Soo-jin Kim"
$ws.Cells.Item(9,7).Value = "'"
$ws.Cells.Item(9,8).Value = "'"
$ws.Cells.Item(9,9).Value = "'"
$ws.Cells.Item(9,10).Value = "Avant-Première en Compétition"

# Row 10
$ws.Cells.Item(10,1).Value = "'2026-02-08"
$ws.Cells.Item(10,2).Value = "'11:00"
$ws.Cells.Item(10,3).Value = "Les fleurs du manguier"
$ws.Cells.Item(10,4).Value = "VO"
$ws.Cells.Item(10,5).Value = "'"
$ws.Cells.Item(10,6).Value = "Akio Fujimoto"
$ws.Cells.Item(10,7).Value = "'"
$ws.Cells.Item(10,8).Value = "'"
$ws.Cells.Item(10,9).Value = "'"
$ws.Cells.Item(10,10).Value = "Avant-Première en Compétition"

# Row 11
$ws.Cells.Item(11,1).Value = "'2026-02-08"
$ws.Cells.Item(11,2).Value = "'15:00"
$ws.Cells.Item(11,3).Value = "Sauvage"
$ws.Cells.Item(11,4).Value = "VF"
$ws.Cells.Item(11,5).Value = "'"
$ws.Cells.Item(11,6).Value = "Camille Ponsin"
$ws.Cells.Item(11,7).Value = "'"
$ws.Cells.Item(11,8).Value = "'"
$ws.Cells.Item(11,9).Value = "'"
$ws.Cells.Item(11,10).Value = "Avant-Première en Compétition"

# Row 12
$ws.Cells.Item(12,1).Value = "'2026-02-08"
$ws.Cells.Item(12,2).Value = "'17:00"
$ws.Cells.Item(12,3).Value = "I swear"
$ws.Cells.Item(12,4).Value = "VO"
$ws.Cells.Item(12,5).Value = "'"
$ws.Cells.Item(12,6).Value = "Kirk Jones"
$ws.Cells.Item(12,7).Value = "'"
$ws.Cells.Item(12,8).Value = "'"
$ws.Cells.Item(12,9).Value = "'"
$ws.Cells.Item(12,10).Value = "Avant-Première en Compétition"

# Row 13
$ws.Cells.Item(13,1).Value = "'2026-02-10"
$ws.Cells.Item(13,2).Value = "'21:00"
$ws.Cells.Item(13,3).Value = "Ma frère"
$ws.Cells.Item(13,4).Value = "VF"
$ws.Cells.Item(13,5).Value = "'"
$ws.Cells.Item(13,6).Value = "Lise Akoka, Romane Gueret"
$ws.Cells.Item(13,7).Value = "Prix du jury et Prix du jury presse Pauillac 2025"
$ws.Cells.Item(13,8).Value = "'"
$ws.Cells.Item(13,9).Value = "'"
$ws.Cells.Item(13,10).Value = "'"

# Row 14
$ws.Cells.Item(14,1).Value = "'2026-02-11"
$ws.Cells.Item(14,2).Value = "'21:00"
$ws.Cells.Item(14,3).Value = "Tatouage"
$ws.Cells.Item(14,4).Value = "VO"
$ws.Cells.Item(14,5).Value = "'"
$ws.Cells.Item(14,6).Value = "Yasuzo Masumura"
$ws.Cells.Item(14,7).Value = "Soirée avec Pup En Vol + invité (à préciser) - Partenariat ADRC"
$ws.Cells.Item(14,8).Value = "'"
$ws.Cells.Item(14,9).Value = "'"
$ws.Cells.Item(14,10).Value = "'"

# Row 15
$ws.Cells.Item(15,1).Value = "'2026-02-12"
$ws.Cells.Item(15,2).Value = "'21:00"
$ws.Cells.Item(15,3).Value = "A pied d'œuvre"
$ws.Cells.Item(15,4).Value = "VF"
$ws.Cells.Item(15,5).Value = "'"
$ws.Cells.Item(15,6).Value = "Valérie Donzelli"
$ws.Cells.Item(15,7).Value = "Partenariat ADRC"
$ws.Cells.Item(15,8).Value = "'"
$ws.Cells.Item(15,9).Value = "'"
$ws.Cells.Item(15,10).Value = "'"

# Row 16
$ws.Cells.Item(16,1).Value = "'2026-02-13"
$ws.Cells.Item(16,2).Value = "'19:00"
$ws.Cells.Item(16,3).Value = "Alter ego"
$ws.Cells.Item(16,4).Value = "VF"
$ws.Cells.Item(16,5).Value = "'"
$ws.Cells.Item(16,6).Value = "Nicolas Charlet, Bruno Lavaine"
$ws.Cells.Item(16,7).Value = "'"
$ws.Cells.Item(16,8).Value = "'"
$ws.Cells.Item(16,9).Value = "'"
$ws.Cells.Item(16,10).Value = "'"

# Row 17
$ws.Cells.Item(17,1).Value = "'2026-02-13"
$ws.Cells.Item(17,2).Value = "'21:15"
$ws.Cells.Item(17,3).Value = "Le mystérieux regard du flamant rose"
$ws.Cells.Item(17,4).Value = "VO"
$ws.Cells.Item(17,5).Value = "'"
$ws.Cells.Item(17,6).Value = "Diego Cespedes"
$ws.Cells.Item(17,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(17,8).Value = "'"
$ws.Cells.Item(17,9).Value = "'"
$ws.Cells.Item(17,10).Value = "'"

# Row 18
$ws.Cells.Item(18,1).Value = "'2026-02-14"
$ws.Cells.Item(18,2).Value = "'14:00"
$ws.Cells.Item(18,3).Value = "Un jour avec mon père"
$ws.Cells.Item(18,4).Value = "VO"
$ws.Cells.Item(18,5).Value = "'"
$ws.Cells.Item(18,6).Value = "Akinola Davies"
$ws.Cells.Item(18,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(18,8).Value = "'"
$ws.Cells.Item(18,9).Value = "'"
$ws.Cells.Item(18,10).Value = "Séance avancée pour jurys - pas dans prog"

# Row 19
$ws.Cells.Item(19,1).Value = "'2026-02-14"
$ws.Cells.Item(19,2).Value = "'16:00"
$ws.Cells.Item(19,3).Value = "Maspalomas"
$ws.Cells.Item(19,4).Value = "VO"
$ws.Cells.Item(19,5).Value = "'"
$ws.Cells.Item(19,6).Value = "Aitor Arregi, José Mari Goenaga"
$ws.Cells.Item(19,7).Value = "AP-COMP  Prix Cinema Europa Les Arcs 2025, Prix d'interprétation San Sebastian 2025"
$ws.Cells.Item(19,8).Value = "'"
$ws.Cells.Item(19,9).Value = "'"
$ws.Cells.Item(19,10).Value = "Partenariat Festival DIAM"

# Row 20
$ws.Cells.Item(20,1).Value = "'2026-02-14"
$ws.Cells.Item(20,2).Value = "'18:15"
$ws.Cells.Item(20,3).Value = "La danse des renards"
$ws.Cells.Item(20,4).Value = "VF"
$ws.Cells.Item(20,5).Value = "'"
$ws.Cells.Item(20,6).Value = "Valéry Carnoy"
$ws.Cells.Item(20,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(20,8).Value = "'"
$ws.Cells.Item(20,9).Value = "'"
$ws.Cells.Item(20,10).Value = "Repas partagé entres les 2 films 19h45"

# Row 21
$ws.Cells.Item(21,1).Value = "'2026-02-14"
$ws.Cells.Item(21,2).Value = "'21:15"
$ws.Cells.Item(21,3).Value = "Le garçon qui faisait danser les  collines"
$ws.Cells.Item(21,4).Value = "VO"
$ws.Cells.Item(21,5).Value = "'"
$ws.Cells.Item(21,6).Value = "Georgi M Unkovski"
$ws.Cells.Item(21,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(21,8).Value = "'"
$ws.Cells.Item(21,9).Value = "'"
$ws.Cells.Item(21,10).Value = "'"

# Row 22
$ws.Cells.Item(22,1).Value = "'2026-02-15"
$ws.Cells.Item(22,2).Value = "'11:00"
$ws.Cells.Item(22,3).Value = "Affection, affection"
$ws.Cells.Item(22,4).Value = "VF"
$ws.Cells.Item(22,5).Value = "'"
$ws.Cells.Item(22,6).Value = "Alexia Walther, Maxime Matray"
$ws.Cells.Item(22,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(22,8).Value = "'"
$ws.Cells.Item(22,9).Value = "'"
$ws.Cells.Item(22,10).Value = "'"

# Row 23
$ws.Cells.Item(23,1).Value = "'2026-02-15"
$ws.Cells.Item(23,2).Value = "'15:00"
$ws.Cells.Item(23,3).Value = "Un jour avec mon père"
$ws.Cells.Item(23,4).Value = "VO"
$ws.Cells.Item(23,5).Value = "'"
$ws.Cells.Item(23,6).Value = "Akinola Davies"
$ws.Cells.Item(23,7).Value = "Avant-Première en Compétition"
$ws.Cells.Item(23,8).Value = "'"
$ws.Cells.Item(23,9).Value = "'"
$ws.Cells.Item(23,10).Value = "'"

# Row 24
$ws.Cells.Item(24,1).Value = "'2026-02-15"
$ws.Cells.Item(24,2).Value = "'17:00"
$ws.Cells.Item(24,3).Value = "La maison des femmes"
$ws.Cells.Item(24,4).Value = "VF"
$ws.Cells.Item(24,5).Value = "'"
$ws.Cells.Item(24,6).Value = "Melisa Godet"
$ws.Cells.Item(24,7).Value = "Avant-Première"
$ws.Cells.Item(24,8).Value = "'"
$ws.Cells.Item(24,9).Value = "'"
$ws.Cells.Item(24,10).Value = "'"

# Row 25
$ws.Cells.Item(25,1).Value = "'2026-02-17"
$ws.Cells.Item(25,2).Value = "'09:30"
$ws.Cells.Item(25,3).Value = "Piro Piro"
$ws.Cells.Item(25,4).Value = "VF"
$ws.Cells.Item(25,5).Value = "'"
$ws.Cells.Item(25,6).Value = "Sung-ah Min"
$ws.Cells.Item(25,7).Value = "SCOL"
$ws.Cells.Item(25,8).Value = "'"
$ws.Cells.Item(25,9).Value = "'"
$ws.Cells.Item(25,10).Value = "Ec Labastide Clermont"

# Row 26
$ws.Cells.Item(26,1).Value = "'2026-02-17"
$ws.Cells.Item(26,2).Value = "'14:30"
$ws.Cells.Item(26,3).Value = "Le mécano de la générale"
$ws.Cells.Item(26,4).Value = "VF"
$ws.Cells.Item(26,5).Value = "'"
$ws.Cells.Item(26,6).Value = "Clyde Bruckman, Buster Keaton"
$ws.Cells.Item(26,7).Value = "SCOL"
$ws.Cells.Item(26,8).Value = "'"
$ws.Cells.Item(26,9).Value = "'"
$ws.Cells.Item(26,10).Value = "Coll Mandela Noé"

# Row 27
$ws.Cells.Item(27,1).Value = "'2026-02-17"
$ws.Cells.Item(27,2).Value = "'21:00"
$ws.Cells.Item(27,3).Value = "L'affaire Bojarski"
$ws.Cells.Item(27,4).Value = "VF"
$ws.Cells.Item(27,5).Value = "CM2"
$ws.Cells.Item(27,6).Value = "Jean-Paul Salomé"
$ws.Cells.Item(27,7).Value = "'"
$ws.Cells.Item(27,8).Value = "'"
$ws.Cells.Item(27,9).Value = "'"
$ws.Cells.Item(27,10).Value = "'"

